# Update countries & provincias Spain
#
# The source workbook ("Pais" sheet) lists countries with COVID-style
# stats (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) in columns A-H, one country per
# row, ordered by total cases (column B) descending.
#
# This update:
#   1. Refreshes the running totals for a number of countries.
#   2. Because Ecuador's totals grew past Colombia's, Ecuador now
#      ranks above Colombia -> the two countries swap rows (29/30).
#   3. Because Paraguay's totals grew past Madagascar's, Paraguay now
#      ranks above Madagascar -> the two countries swap rows (118/119).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 2101143
$ws.Range("C4").Value = 11442
$ws.Range("D4").Value = 819375
$ws.Range("E4").Value = 1165371
$ws.Range("G4").Value = 363
$ws.Range("H4").Value = 116397

# --- India (row 7) ---
$ws.Range("B7").Value = 309389
$ws.Range("C7").Value = 11106
$ws.Range("D7").Value = 154131
$ws.Range("E7").Value = 146368
$ws.Range("G7").Value = 389
$ws.Range("H7").Value = 8890

# --- Alemania (row 12) ---
$ws.Range("B12").Value = 187010
$ws.Range("C12").Value = 215
$ws.Range("E12").Value = 6556
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 8854

# --- Turquia (row 14) ---
$ws.Range("B14").Value = 175218
$ws.Range("C14").Value = 1195
$ws.Range("D14").Value = 149102
$ws.Range("E14").Value = 21338
$ws.Range("G14").Value = 15
$ws.Range("H14").Value = 4778

# --- Colombia / Ecuador swap (rows 29-30) ---
# Row 29 becomes Ecuador (updated stats), row 30 becomes Colombia
# (carrying forward its previous, unchanged stats).
$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 45778
$ws.Range("C29").Value = 1338
$ws.Range("D29").Value = 22679
$ws.Range("E29").Value = 19271
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 108
$ws.Range("H29").Value = 3828

$ws.Range("A30").Value = "Colombia"
$ws.Range("B30").Value = 45212
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 17790
$ws.Range("E30").Value = 25934
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 1488

# --- Israel (row 47) ---
$ws.Range("B47").Value = 18795
$ws.Range("C47").Value = 226
$ws.Range("E47").Value = 3207

# --- Kazajistan (row 56) ---
$ws.Range("D56").Value = 8829
$ws.Range("E56").Value = 4973
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 70

# --- Marruecos (row 67) ---
$ws.Range("B67").Value = 8610
$ws.Range("C67").Value = 73
$ws.Range("D67").Value = 7618
$ws.Range("E67").Value = 780

# --- Maldivas (row 101) ---
$ws.Range("B101").Value = 2003
$ws.Range("C101").Value = 27
$ws.Range("D101").Value = 1193
$ws.Range("E101").Value = 802

# --- Madagascar / Paraguay swap (rows 118-119) ---
# Row 118 becomes Paraguay (updated stats), row 119 becomes
# Madagascar (carrying forward its previous, unchanged stats).
$ws.Range("A118").Value = "Paraguay"
$ws.Range("B118").Value = 1254
$ws.Range("C118").Value = 24
$ws.Range("D118").Value = 633
$ws.Range("E118").Value = 610
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 11

$ws.Range("A119").Value = "Madagascar"
$ws.Range("B119").Value = 1240
$ws.Range("C119").Value = 37
$ws.Range("D119").Value = 344
$ws.Range("E119").Value = 886
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 10

# --- Suazilandia (row 150) ---
$ws.Range("B150").Value = 472
$ws.Range("C150").Value = 23
$ws.Range("D150").Value = 246
$ws.Range("E150").Value = 223
